{"js": "// The edit replaces the date in the title paragraph and the 25\n// multiplication-problem strings found in the non-empty cells of the\n// single table, each with a new value, in document order. Every text\n// item in the document is changed (old values are all unique), so we\n// can safely walk the body in order and apply the new values in turn.\n\nconst oldToNew = new Map([\n  [\"2024-07-16 Tuesday\", \"2024-07-17 Wednesday\"],\n  [\"603\u00d79=\", \"305\u00d79=\"],\n  [\"870\u00d75=\", \"868\u00d79=\"],\n  [\"597\u00d79=\", \"380\u00d77=\"],\n  [\"930\u00d75=\", \"559\u00d79=\"],\n  [\"246\u00d75=\", \"451\u00d75=\"],\n  [\"182\u00d74=\", \"793\u00d75=\"],\n  [\"618\u00d76=\", \"748\u00d72=\"],\n  [\"463\u00d72=\", \"758\u00d74=\"],\n  [\"785\u00d75=\", \"984\u00d77=\"],\n  [\"966\u00d76=\", \"296\u00d73=\"],\n  [\"798\u00d74=\", \"822\u00d75=\"],\n  [\"432\u00d78=\", \"220\u00d74=\"],\n  [\"815\u00d75=\", \"425\u00d77=\"],\n  [\"157\u00d76=\", \"969\u00d72=\"],\n  [\"723\u00d72=\", \"900\u00d78=\"],\n  [\"489\u00d72=\", \"118\u00d72=\"],\n  [\"556\u00d74=\", \"315\u00d74=\"],\n  [\"454\u00d79=\", \"197\u00d72=\"],\n  [\"557\u00d79=\", \"685\u00d73=\"],\n  [\"989\u00d73=\", \"566\u00d79=\"],\n  [\"479\u00d72=\", \"670\u00d72=\"],\n  [\"486\u00d74=\", \"672\u00d74=\"],\n  [\"123\u00d72=\", \"163\u00d75=\"],\n  [\"663\u00d76=\", \"751\u00d78=\"],\n  [\"260\u00d76=\", \"833\u00d75=\"],\n]);\n\n// 1) The title paragraph holding the date (the first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const newText = oldToNew.get(p.text);\n  if (newText !== undefined) {\n    p.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// 2) The problem strings living in the table cells.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  table.load(\"values,rowCount\");\n  await context.sync();\n\n  for (let r = 0; r < table.rowCount; r++) {\n    const rowValues = table.values[r];\n    for (let c = 0; c < rowValues.length; c++) {\n      const cellText = rowValues[c];\n      if (cellText === \"\") {\n        continue;\n      }\n      const newText = oldToNew.get(cellText);\n      if (newText !== undefined) {\n        table.getCell(r, c).value = newText;\n      }\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# The edit replaces the date in the title paragraph and the 25\n# multiplication-problem strings found in the table cells, each with a\n# new value. Every old value below is unique within the document, so a\n# simple whole-document Find/Replace for each pair reproduces the diff\n# exactly while leaving all other formatting untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-16 Tuesday\", \"2024-07-17 Wednesday\"),\n    @(\"603\u00d79=\", \"305\u00d79=\"),\n    @(\"870\u00d75=\", \"868\u00d79=\"),\n    @(\"597\u00d79=\", \"380\u00d77=\"),\n    @(\"930\u00d75=\", \"559\u00d79=\"),\n    @(\"246\u00d75=\", \"451\u00d75=\"),\n    @(\"182\u00d74=\", \"793\u00d75=\"),\n    @(\"618\u00d76=\", \"748\u00d72=\"),\n    @(\"463\u00d72=\", \"758\u00d74=\"),\n    @(\"785\u00d75=\", \"984\u00d77=\"),\n    @(\"966\u00d76=\", \"296\u00d73=\"),\n    @(\"798\u00d74=\", \"822\u00d75=\"),\n    @(\"432\u00d78=\", \"220\u00d74=\"),\n    @(\"815\u00d75=\", \"425\u00d77=\"),\n    @(\"157\u00d76=\", \"969\u00d72=\"),\n    @(\"723\u00d72=\", \"900\u00d78=\"),\n    @(\"489\u00d72=\", \"118\u00d72=\"),\n    @(\"556\u00d74=\", \"315\u00d74=\"),\n    @(\"454\u00d79=\", \"197\u00d72=\"),\n    @(\"557\u00d79=\", \"685\u00d73=\"),\n    @(\"989\u00d73=\", \"566\u00d79=\"),\n    @(\"479\u00d72=\", \"670\u00d72=\"),\n    @(\"486\u00d74=\", \"672\u00d74=\"),\n    @(\"123\u00d72=\", \"163\u00d75=\"),\n    @(\"663\u00d76=\", \"751\u00d78=\"),\n    @(\"260\u00d76=\", \"833\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = \"wdFindContinue\"\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, \"wdFindContinue\", $false, $new, \"wdReplaceAll\") | Out-Null\n}\n"}
